$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = -13
$ws.Range("F11").Value = -6
$ws.Range("F12").Value = -4
